$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "hub_ids" column (G) entirely - specialization_ids/joining_year
# shift left into G/H.
$ws.Columns.Item(7).Delete()

# Corrected mobile numbers, now stored as text instead of numbers.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "9909990132"

$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "7990357110"

$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "9753555363"

$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "8790818948"

# Move/restore the active selection to G1 (where hub_ids used to be).
$ws.Range("G1").Select() | Out-Null
